# Add a new "ShowBlip" boolean flag to the property Flags bit-field.
# - Updates the explanatory text on the "Field Explanation" sheet.
# - Updates the numeric Flags values on the "Properties Table" sheet
#   (appending a new low-order digit of 1 for every row, since the new
#   flag defaults to "on"/shown).
# - Leaves the active sheet/selection on "Properties Table".

$wb = $excel.ActiveWorkbook

# --- 1. "Field Explanation" sheet: document the new flag -------------------
$wsField = $wb.Worksheets.Item("Field Explanation")
$wsField.Range("C9").Value = "Used as 'bulk storage' for the boolean variables [Ownable][Owned][ContextMission][ShowBlip]"

# --- 2. "Properties Table" sheet: extend the Flags column -------------------
$wsProps = $wb.Worksheets.Item("Properties Table")

$wsProps.Range("C2").Value = 1011
$wsProps.Range("C3").Value = 1002
$wsProps.Range("C4").Value = 1001
$wsProps.Range("C5").Value = 1011
$wsProps.Range("C6").Value = 1001
$wsProps.Range("C7").Value = 1001
$wsProps.Range("C8").Value = 1001
$wsProps.Range("C9").Value = 1001
$wsProps.Range("C10").Value = 1001
$wsProps.Range("C11").Value = 1001

# --- 3. Update view/selection state -----------------------------------------
# "Field Explanation" no longer the active tab; its lingering selection moves
# to C9 (where the edit was made).
$wsField.Activate()
$wsField.Range("C9").Select()

# "Properties Table" becomes the active tab, with the selection left on B13.
$wsProps.Activate()
$wsProps.Range("B13").Select()
